$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data-vintage note in the source line (avril 2023 -> octobre 2023)
$sourceCell = $ws.Range("A104")
$sourceCell.Value = $sourceCell.Value2 -replace "avril 2023", "octobre 2023"

# Refresh the computed GDP-expenditure-decomposition figures (rows recalculated
# with the October 2023 WEO vintage)
$ws.Range("C13").Value = 78.5051630092483 ; $ws.Range("D13").Value = 18.5774860555573 ; $ws.Range("E13").Value = 0.7775280234092 ; $ws.Range("F13").Value = 2.13982291178524 ; $ws.Range("G13").Value = 33.3687189901663 ; $ws.Range("H13").Value = 35.5085419019515
$ws.Range("C23").Value = 79.2174750553685 ; $ws.Range("D23").Value = 17.6646574191516 ; $ws.Range("E23").Value = 0.05303801960893 ; $ws.Range("F23").Value = 3.06482950587096 ; $ws.Range("G23").Value = 35.1284306370612 ; $ws.Range("H23").Value = 38.1932601429322
$ws.Range("C38").Value = 82.608107524147 ; $ws.Range("D38").Value = 26.51321769747 ; $ws.Range("E38").Value = -0.3310465520867 ; $ws.Range("F38").Value = -8.7902786695303 ; $ws.Range("G38").Value = 24.3340102339887 ; $ws.Range("H38").Value = 15.5437315644585
$ws.Range("C45").Value = 83.0240191458582 ; $ws.Range("D45").Value = 18.921899364408 ; $ws.Range("E45").Value = 2.63773116683632 ; $ws.Range("F45").Value = -4.5836496771024 ; $ws.Range("G45").Value = 29.0379864718626 ; $ws.Range("H45").Value = 24.4543367947601
$ws.Range("C61").Value = 81.3109691816165 ; $ws.Range("D61").Value = 21.6173313623917 ; $ws.Range("E61").Value = 0.84944230941719 ; $ws.Range("F61").Value = -3.7777428534254 ; $ws.Range("G61").Value = 23.3913634702131 ; $ws.Range("H61").Value = 19.6136206167877
$ws.Range("C62").Value = 81.4850898488542 ; $ws.Range("D62").Value = 20.612459349262 ; $ws.Range("E62").Value = 1.25886325432601 ; $ws.Range("F62").Value = -3.3564124524421 ; $ws.Range("G62").Value = 27.98218780978 ; $ws.Range("H62").Value = 24.6257753573379
$ws.Range("C63").Value = 70.7800615949826 ; $ws.Range("D63").Value = 27.5906993476265 ; $ws.Range("E63").Value = 1.30750430955534 ; $ws.Range("F63").Value = 0.32173474783554 ; $ws.Range("G63").Value = 30.4868915615395 ; $ws.Range("H63").Value = 30.7593272443428
$ws.Range("C64").Value = 81.3312930878922 ; $ws.Range("D64").Value = 20.2839356963998 ; $ws.Range("E64").Value = 0.31616236735951 ; $ws.Range("F64").Value = -1.9313911516515 ; $ws.Range("G64").Value = 30.7633814989786 ; $ws.Range("H64").Value = 28.8319903473271
$ws.Range("C65").Value = 62.3336815559125 ; $ws.Range("D65").Value = 36.318879054108 ; $ws.Range("E65").Value = 1.2589130121668 ; $ws.Range("F65").Value = 0.0885263778127 ; $ws.Range("G65").Value = 23.6596735498598 ; $ws.Range("H65").Value = 23.6813895489389
$ws.Range("C66").Value = 71.3085901856593 ; $ws.Range("D66").Value = 27.2461697243249 ; $ws.Range("E66").Value = 1.3051028036686 ; $ws.Range("F66").Value = 0.14013728634719 ; $ws.Range("G66").Value = 30.3621985703611 ; $ws.Range("H66").Value = 30.4539773877773
$ws.Range("C67").Value = 88.5477485006933 ; $ws.Range("D67").Value = 16.6857375053157 ; $ws.Range("E67").Value = 2.38107712094283 ; $ws.Range("F67").Value = -7.6145631269518 ; $ws.Range("G67").Value = 26.4611379451883 ; $ws.Range("H67").Value = 18.8465748182365
$ws.Range("C68").Value = 85.2290807314893 ; $ws.Range("D68").Value = 18.520924284648 ; $ws.Range("E68").Value = 2.10251580139608 ; $ws.Range("F68").Value = -5.8525208175334 ; $ws.Range("G68").Value = 27.1929748428028 ; $ws.Range("H68").Value = 21.3404540252694
$ws.Range("C69").Value = 82.531416543219 ; $ws.Range("D69").Value = 25.3152188682017 ; $ws.Range("E69").Value = -0.4242350086144 ; $ws.Range("F69").Value = -7.4224004028063 ; $ws.Range("G69").Value = 26.3250290444988 ; $ws.Range("H69").Value = 18.9026286416924
$ws.Range("C70").Value = 72.2728779727057 ; $ws.Range("D70").Value = 20.2338752756904 ; $ws.Range("E70").Value = 0.01089896364686 ; $ws.Range("F70").Value = 7.48234778795711 ; $ws.Range("G70").Value = 30.5334534631329 ; $ws.Range("H70").Value = 38.01580125109
$ws.Range("C71").Value = 81.3109691816165 ; $ws.Range("D71").Value = 21.6173313623917 ; $ws.Range("E71").Value = 0.84944230941719 ; $ws.Range("F71").Value = -3.7777428534254 ; $ws.Range("G71").Value = 23.3913634702131 ; $ws.Range("H71").Value = 19.6136206167877
$ws.Range("C72").Value = 85.962199116529 ; $ws.Range("D72").Value = 23.6485834980376 ; $ws.Range("E72").Value = 0.10227472358688 ; $ws.Range("F72").Value = -9.7130573381534 ; $ws.Range("G72").Value = 21.7599072369958 ; $ws.Range("H72").Value = 12.0468498988424
$ws.Range("C73").Value = 78.2535287273749 ; $ws.Range("D73").Value = 20.7152542199366 ; $ws.Range("E73").Value = 0.3594157416111 ; $ws.Range("F73").Value = 0.67180131107747 ; $ws.Range("G73").Value = 33.4882359480136 ; $ws.Range("H73").Value = 34.160037259091
$ws.Range("C74").Value = 73.0124498675532 ; $ws.Range("D74").Value = 27.9419885488458 ; $ws.Range("E74").Value = 0.33641836968079 ; $ws.Range("F74").Value = -1.2908567860798 ; $ws.Range("G74").Value = 39.6054279933623 ; $ws.Range("H74").Value = 38.3145712072826
$ws.Range("C75").Value = 64.1603650160263 ; $ws.Range("D75").Value = 29.8973578636498 ; $ws.Range("E75").Value = 0.02178612424328 ; $ws.Range("F75").Value = 5.92049099608062 ; $ws.Range("G75").Value = 33.516792465615 ; $ws.Range("H75").Value = 39.4372834616956
$ws.Range("C76").Value = 68.1033318775114 ; $ws.Range("D76").Value = 26.9959767662255 ; $ws.Range("E76").Value = 2.0641224566128 ; $ws.Range("F76").Value = 2.83656889965034 ; $ws.Range("G76").Value = 50.4357891900293 ; $ws.Range("H76").Value = 53.260400910543
$ws.Range("C77").Value = 80.6310570663695 ; $ws.Range("D77").Value = 19.6812107050154 ; $ws.Range("E77").Value = 0.26863656005942 ; $ws.Range("F77").Value = -0.5809043314443 ; $ws.Range("G77").Value = 22.9814324458976 ; $ws.Range("H77").Value = 22.4005281144533
$ws.Range("C78").Value = 74.2583297490871 ; $ws.Range("D78").Value = 22.5806667216793 ; $ws.Range("E78").Value = 1.97878214318177 ; $ws.Range("F78").Value = 1.18222138605182 ; $ws.Range("G78").Value = 54.5490089744429 ; $ws.Range("H78").Value = 55.7312303604948
$ws.Range("C79").Value = 77.3155980445935 ; $ws.Range("D79").Value = 22.715629722513 ; $ws.Range("E79").Value = 1.32791258245413 ; $ws.Range("F79").Value = -1.3591403495606 ; $ws.Range("G79").Value = 34.7739631889701 ; $ws.Range("H79").Value = 33.414823206107
$ws.Range("C80").Value = 60.2010841988602 ; $ws.Range("D80").Value = 29.5496423231308 ; $ws.Range("E80").Value = -0.7354713018934 ; $ws.Range("F80").Value = 10.9847447799024 ; $ws.Range("G80").Value = 26.7068380207434 ; $ws.Range("H80").Value = 37.6915828006457
$ws.Range("C81").Value = 61.6808808155568 ; $ws.Range("D81").Value = 22.3327252369851 ; $ws.Range("E81").Value = 3.79983986854925 ; $ws.Range("F81").Value = 12.1865540789089 ; $ws.Range("G81").Value = 25.3461903609255 ; $ws.Range("H81").Value = 37.1885680736754
$ws.Range("C82").Value = 84.5719466042157 ; $ws.Range("D82").Value = 19.3162839915564 ; $ws.Range("E82").Value = 1.54810510092919 ; $ws.Range("F82").Value = -5.4363356967013 ; $ws.Range("G82").Value = 28.1671540316687 ; $ws.Range("H82").Value = 22.7308183349674
$ws.Range("C83").Value = 71.5328766658131 ; $ws.Range("D83").Value = 28.0257145702581 ; $ws.Range("E83").Value = 1.1013024749221 ; $ws.Range("F83").Value = -0.6598937109932 ; $ws.Range("G83").Value = 30.9162467720621 ; $ws.Range("H83").Value = 30.2223522399506
$ws.Range("C84").Value = 85.1771917562614 ; $ws.Range("D84").Value = 24.3970454768094 ; $ws.Range("E84").Value = 0.28280812965316 ; $ws.Range("F84").Value = -9.8570453627239 ; $ws.Range("G84").Value = 33.0879080775183 ; $ws.Range("H84").Value = 23.2308627147943
$ws.Range("C86").Value = 80.6277565894622 ; $ws.Range("D86").Value = 20.8979172418502 ; $ws.Range("E86").Value = 1.53127633348599 ; $ws.Range("F86").Value = -3.0569501647984 ; $ws.Range("G86").Value = 25.6385509639328 ; $ws.Range("H86").Value = 22.5816007991344
$ws.Range("C87").Value = 73.7667471034045 ; $ws.Range("D87").Value = 28.7321530533509 ; $ws.Range("E87").Value = 1.89553179182588 ; $ws.Range("F87").Value = -4.3944319485812 ; $ws.Range("G87").Value = 27.444438828795 ; $ws.Range("H87").Value = 22.6423511295252
$ws.Range("C88").Value = 81.8402568866825 ; $ws.Range("D88").Value = 15.241155485426 ; $ws.Range("E88").Value = 0.94478644766041 ; $ws.Range("F88").Value = 1.97380118023119 ; $ws.Range("G88").Value = 33.8555740775008 ; $ws.Range("H88").Value = 35.829375257732
$ws.Range("C89").Value = 62.1916136194635 ; $ws.Range("D89").Value = 34.0619383581071 ; $ws.Range("E89").Value = 1.13757617912919 ; $ws.Range("F89").Value = 2.60887184330023 ; $ws.Range("G89").Value = 24.7932386976649 ; $ws.Range("H89").Value = 27.3272591362631
$ws.Range("C90").Value = 75.8926279013926 ; $ws.Range("D90").Value = 22.5679837489585 ; $ws.Range("E90").Value = 1.24070343644199 ; $ws.Range("F90").Value = 0.29868491320697 ; $ws.Range("G90").Value = 35.6492941107351 ; $ws.Range("H90").Value = 35.9479794005185
$ws.Range("C91").Value = 78.6759516837579 ; $ws.Range("D91").Value = 27.2920476078795 ; $ws.Range("E91").Value = -0.1846825520048 ; $ws.Range("F91").Value = -5.7833167396326 ; $ws.Range("G91").Value = 31.2674734361747 ; $ws.Range("H91").Value = 25.4841566965421
$ws.Range("C92").Value = 80.9583754183292 ; $ws.Range("D92").Value = 30.6921029362763 ; $ws.Range("E92").Value = 0.77887839135742 ; $ws.Range("F92").Value = -12.429356745963 ; $ws.Range("G92").Value = 29.6808953926222 ; $ws.Range("H92").Value = 17.2515386466593
$ws.Range("C93").Value = 89.5610474142206 ; $ws.Range("D93").Value = 22.1913750948632 ; $ws.Range("E93").Value = -0.0530574162454 ; $ws.Range("F93").Value = -11.699365092838 ; $ws.Range("G93").Value = 62.8749532252898 ; $ws.Range("H93").Value = 51.1755881324516
$ws.Range("C94").Value = 56.3217364865289 ; $ws.Range("D94").Value = 23.1726500054126 ; $ws.Range("E94").Value = 0.99181163759275 ; $ws.Range("F94").Value = 19.5138018704658 ; $ws.Range("G94").Value = 110.131841241474 ; $ws.Range("H94").Value = 129.64564311194
$ws.Range("C95").Value = 83.6207304743834 ; $ws.Range("D95").Value = 24.7119970188535 ; $ws.Range("E95").Value = 0.38851897979537 ; $ws.Range("F95").Value = -8.7212464730323 ; $ws.Range("G95").Value = 28.6069113992903 ; $ws.Range("H95").Value = 19.8856649262581
$ws.Range("C96").Value = 74.1156976229987 ; $ws.Range("D96").Value = 24.4175957286443 ; $ws.Range("E96").Value = 3.55715731203879 ; $ws.Range("F96").Value = -2.0904506636818 ; $ws.Range("G96").Value = 37.8979120561488 ; $ws.Range("H96").Value = 35.8074613924671
$ws.Range("C97").Value = 79.691397980803 ; $ws.Range("D97").Value = 24.5291070960319 ; $ws.Range("E97").Value = -0.0866356206933 ; $ws.Range("F97").Value = -4.1338694561416 ; $ws.Range("G97").Value = 23.3084189956619 ; $ws.Range("H97").Value = 19.1745495395203
$ws.Range("C98").Value = 78.6042774702587 ; $ws.Range("D98").Value = 23.151717389688 ; $ws.Range("E98").Value = 4.61922461544348 ; $ws.Range("F98").Value = -6.3752194753902 ; $ws.Range("G98").Value = 23.914261042786 ; $ws.Range("H98").Value = 17.5390415673958
